# Auto-generated Excel COM-interop script to apply crypto price/volume updates
# per commit "Updated cryptos list on Sat Jan 13 22:37:25 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($ref, $text) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

$updates = @(
    @('D2', '42.919.73'),
    @('E2', '  +1.62%  '),
    @('D3', '2.570.82'),
    @('E3', '  +2.84%  '),
    @('D4', '0.999'),
    @('E4', '  +0.37%  '),
    @('D5', '302.65'),
    @('E5', '  +3.16%  '),
    @('D6', '96.70'),
    @('E6', '  +4.88%  '),
    @('D7', '0.576'),
    @('E7', '  +1.66%  '),
    @('E8', '  +0.11%  '),
    @('D9', '0.549'),
    @('E9', '  +1.39%  '),
    @('D10', '36.74'),
    @('E10', '  +2.64%  '),
    @('D11', '0.0807'),
    @('E11', '  +1.94%  '),
    @('D12', '7.73'),
    @('E12', '  +2.01%  '),
    @('E13', '  +7.58%  '),
    @('D14', '2.581.72'),
    @('E14', '  +3.99%  '),
    @('D15', '0.886'),
    @('E15', '  +3.68%  '),
    @('E16', '  +3.00%  '),
    @('D17', '42.941.16'),
    @('E17', '  +2.06%  '),
    @('D18', '13.12'),
    @('E18', '  +7.86%  '),
    @('D19', '0.0₃0993'),
    @('E19', '  +4.34%  '),
    @('E20', '  +3.20%  '),
    @('D21', '72.02'),
    @('E21', '  +0.37%  '),
    @('D22', '254.34'),
    @('E22', '  -0.90%  '),
    @('D23', '2.96'),
    @('E23', '  +3.69%  '),
    @('D24', '2.12'),
    @('E24', '  +1.52%  '),
    @('D25', '28.57'),
    @('E25', '  +0.84%  '),
    @('E26', '  +0.05%  '),
    @('D27', '10.24'),
    @('E27', '  +4.26%  '),
    @('D28', '37.76'),
    @('E28', '  +3.60%  '),
    @('E29', '  -4.07%  '),
    @('D30', '6.07'),
    @('E30', '  +2.42%  '),
    @('D31', '155.63'),
    @('E31', '  +3.82%  '),
    @('E32', '  +1.19%  '),
    @('D33', '2.77'),
    @('E33', '  +2.81%  '),
    @('E34', '  -1.26%  '),
    @('D35', '0.0810'),
    @('E35', '  +2.99%  '),
    @('D36', '18.31'),
    @('E36', '  +10.22%  '),
    @('E37', '  +2.04%  '),
    @('E38', '  +1.72%  '),
    @('D39', '23.62'),
    @('E39', '  +0.45%  '),
    @('D40', '3.42'),
    @('E40', '  +0.66%  '),
    @('B41', 'ApeXProtocol'),
    @('C41', 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'),
    @('D41', '2.07'),
    @('E41', '  +29.89%  '),
    @('B42', 'RenderToken'),
    @('C42', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'),
    @('D42', '3.87'),
    @('E42', '  +2.81%  '),
    @('B43', 'VeChain'),
    @('C43', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @('D43', '0.0310'),
    @('E43', '  +1.87%  '),
    @('D44', '2.064.34'),
    @('E44', '  +3.41%  '),
    @('D45', '0.998'),
    @('E45', '  +0.60%  '),
    @('E46', '  +6.09%  '),
    @('D47', '85.06'),
    @('E47', '  +0.30%  '),
    @('D48', '77.37'),
    @('E48', '  +15.31%  '),
    @('D49', '2.822.82'),
    @('E49', '  +3.55%  '),
    @('D50', '106.04'),
    @('E50', '  +4.45%  '),
    @('D51', '0.191'),
    @('E51', '  +3.50%  ')
)

foreach ($u in $updates) {
    Set-Text $u[0] $u[1]
}
